$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$msg5f = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e025622058dd88a41bf5ae265f9e63e081f58d2f/e2e/5f266658-4f52-4a5f-9ade-256926e29f04.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/343c0a593262a7f0dabe1a0c956319e395f9f858/e2e/5f266658-4f52-4a5f-9ade-256926e29f04.md."
$msg6d = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e025622058dd88a41bf5ae265f9e63e081f58d2f/e2e/6d4abbe6-e83c-4be9-9a23-aa84760ac57b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/343c0a593262a7f0dabe1a0c956319e395f9f858/e2e/6d4abbe6-e83c-4be9-9a23-aa84760ac57b.md."

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-10-25 03:08:29"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-25 03:08:29"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2016-10-25 03:08:16"
$wsZhCn.Range("P2").Value = $msg5f

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-10-25 03:08:16"
$wsZhCn.Range("P3").Value = $msg6d

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("H2").Value = "2016-10-25 03:08:29"
$wsDeDe.Range("P2").Value = $msg5f

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-10-25 03:08:29"
$wsDeDe.Range("P3").Value = $msg6d

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(16).ColumnWidth = 40
